# Weekly update: insert two new "Caigua" price rows (29-10-2022 week,
# serial 44858) ahead of the existing history, pushing the prior rows
# (old 129-142) down to 131-144.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("129:130").Insert()

$ws.Cells.Item(129, 1).Value = 1
$ws.Cells.Item(129, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(129, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(129, 4).Value = 44858
$ws.Cells.Item(129, 5).Value = 15
$ws.Cells.Item(129, 6).Value = 100112036
$ws.Cells.Item(129, 7).Value = "Caigua"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 140
$ws.Cells.Item(129, 11).Value = 7000
$ws.Cells.Item(129, 12).Value = 8000
$ws.Cells.Item(129, 13).Value = 7500
$ws.Cells.Item(129, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(129, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(129, 16).Value = 375
$ws.Cells.Item(129, 17).Value = 20
$ws.Cells.Item(129, 18).Value = "Hortaliza"

$ws.Cells.Item(130, 1).Value = 1
$ws.Cells.Item(130, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(130, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(130, 4).Value = 44858
$ws.Cells.Item(130, 5).Value = 15
$ws.Cells.Item(130, 6).Value = 100112036
$ws.Cells.Item(130, 7).Value = "Caigua"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Segunda"
$ws.Cells.Item(130, 10).Value = 150
$ws.Cells.Item(130, 11).Value = 6000
$ws.Cells.Item(130, 12).Value = 7000
$ws.Cells.Item(130, 13).Value = 6500
$ws.Cells.Item(130, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(130, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(130, 16).Value = 325
$ws.Cells.Item(130, 17).Value = 20
$ws.Cells.Item(130, 18).Value = "Hortaliza"
